$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint1")

# --- Row 3 ("Tache 1") burndown values ---
# H3 changes 2 -> 4; P3:X3 were empty (style 8) and become 0 with style 19
# (matching the style already used by the rest of the row, G3:O3).
$ws.Range("O3").Copy()
$ws.Range("P3:X3").PasteSpecial(-4122)

$ws.Range("H3").Value = 4
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 0
$ws.Range("S3").Value = 0
$ws.Range("T3").Value = 0
$ws.Range("U3").Value = 0
$ws.Range("V3").Value = 0
$ws.Range("W3").Value = 0
$ws.Range("X3").Value = 0

# --- Row 4 ("Tache 2") burndown values ---
# J4:O4 flatten to 16; P4:X4 (already styled) get filled in, X4 = 0
$ws.Range("J4").Value = 16
$ws.Range("K4").Value = 16
$ws.Range("L4").Value = 16
$ws.Range("M4").Value = 16
$ws.Range("N4").Value = 16
$ws.Range("O4").Value = 16
$ws.Range("P4").Value = 16
$ws.Range("Q4").Value = 16
$ws.Range("R4").Value = 16
$ws.Range("S4").Value = 16
$ws.Range("T4").Value = 16
$ws.Range("U4").Value = 16
$ws.Range("V4").Value = 16
$ws.Range("W4").Value = 16
$ws.Range("X4").Value = 0

# --- Row 5 ("Tache 3") burndown values ---
# P5:X5 (already styled) get filled in with 10
$ws.Range("P5").Value = 10
$ws.Range("Q5").Value = 10
$ws.Range("R5").Value = 10
$ws.Range("S5").Value = 10
$ws.Range("T5").Value = 10
$ws.Range("U5").Value = 10
$ws.Range("V5").Value = 10
$ws.Range("W5").Value = 10
$ws.Range("X5").Value = 10

# --- Row 6 ("Tache 4") burndown values ---
# J6:N6 flatten to 10; P6:X6 (already styled) get filled in with 0
$ws.Range("J6").Value = 10
$ws.Range("K6").Value = 10
$ws.Range("L6").Value = 10
$ws.Range("M6").Value = 10
$ws.Range("N6").Value = 10
$ws.Range("P6").Value = 0
$ws.Range("Q6").Value = 0
$ws.Range("R6").Value = 0
$ws.Range("S6").Value = 0
$ws.Range("T6").Value = 0
$ws.Range("U6").Value = 0
$ws.Range("V6").Value = 0
$ws.Range("W6").Value = 0
$ws.Range("X6").Value = 0

# --- Selection moved from AE8 to AB24 ---
$ws.Range("AB24").Select()
